# Update "想去人数" (F column) figures on both the "展览" and "全部类型" sheets.
# Source data refreshed as of commit 456a3b4 (gh-pages generated output).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - row -> new F value
$sheetExhibition = $wb.Worksheets.Item("展览")
$exhibitionUpdates = @{
    2  = 605
    4  = 1280
    6  = 14229
    7  = 16108
    9  = 77
    10 = 48
    18 = 95
    19 = 33
    20 = 1237
    23 = 27
    24 = 6434
    26 = 7
    29 = 5665
    30 = 90
    32 = 160
    33 = 4684
    34 = 11
}
foreach ($row in $exhibitionUpdates.Keys) {
    $sheetExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

# Sheet "全部类型" (all types) - row -> new F value
$sheetAll = $wb.Worksheets.Item("全部类型")
$allUpdates = @{
    2  = 605
    4  = 1280
    6  = 14229
    7  = 16108
    9  = 77
    10 = 48
    18 = 95
    19 = 33
    20 = 1237
    24 = 27
    25 = 6434
    27 = 7
    31 = 5665
    32 = 90
    34 = 160
    35 = 4684
    36 = 11
}
foreach ($row in $allUpdates.Keys) {
    $sheetAll.Range("F$row").Value = $allUpdates[$row]
}
